# Auto-update draw results: append the 2025-11-23 "Pick 4" result as a new
# row at the bottom of the "Results" sheet, matching the existing table
# layout (Date, Game, Phase, Result, InsertedAt) where every value is
# stored as literal text (e.g. the "251123" phase code and the
# "2025-11-23" date are text, not a number/date serial).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Results")

$newRow = $ws.UsedRange.Rows.Count + 1

$values = @("2025-11-23", "Pick 4", "251123", "0-7-5-8", "2025-11-23T21:37:33.669+04:00")

for ($col = 1; $col -le 5; $col++) {
    $cell = $ws.Cells.Item($newRow, $col)

    # Some of the values look like a date ("2025-11-23") or a plain number
    # ("251123") and would otherwise be auto-converted to a date serial /
    # numeric value on assignment. Briefly force text storage so the literal
    # string is kept (matching every other row in this column), then drop
    # the now-unneeded explicit formatting again so the new cells end up
    # styled exactly like the existing ones (no explicit style index).
    $cell.NumberFormat = "@"
    $cell.Value = $values[$col - 1]
    $cell.ClearFormats()
}
